$wb = $excel.ActiveWorkbook

# ALC (sheet1) row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value2 = 2300.4167
$ws.Range("I116").Value2 = 1837.5
$ws.Range("J116").Value2 = 3226.25
$ws.Range("K116").Value2 = 1837.5
$ws.Range("L116").Value2 = 3226.25
$ws.Range("M116").Value2 = 1604.5
$ws.Range("N116").Value2 = -10110.25

# ARM (sheet2) row 3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value2 = 500
$ws.Range("J3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("N3").ClearContents()

# ARM (sheet2) row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3675.74
$ws.Range("I32").Value2 = 3675.74
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 3675.74
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -3388.74
$ws.Range("N32").ClearContents()

# ARM (sheet2) row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1606.0164
$ws.Range("I61").Value2 = 1519.5946
$ws.Range("J61").Value2 = 1739.25
$ws.Range("K61").Value2 = 1519.5946
$ws.Range("L61").Value2 = 1739.25
$ws.Range("M61").Value2 = -1307.5946
$ws.Range("N61").Value2 = -2163.25

# ARM (sheet2) row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 1145.9
$ws.Range("I122").Value2 = 1120.125
$ws.Range("J122").Value2 = 1249
$ws.Range("K122").Value2 = 3360.375
$ws.Range("L122").Value2 = 3747
$ws.Range("M122").Value2 = -910.375
$ws.Range("N122").Value2 = -8647

# ARM (sheet2) row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value2 = 31475.916
$ws.Range("J125").Value2 = 31475.916
$ws.Range("L125").Value2 = 31475.916
$ws.Range("N125").Value2 = -41315.916

# ARM (sheet2) row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 4708.364
$ws.Range("I132").Value2 = 1459.7727
$ws.Range("K132").Value2 = 4379.3181
$ws.Range("M132").Value2 = -1849.3181

# ARM (sheet2) row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 1606.0164
$ws.Range("I136").Value2 = 1519.5946
$ws.Range("J136").Value2 = 1739.25
$ws.Range("K136").Value2 = 4558.783799999999
$ws.Range("L136").Value2 = 5217.75
$ws.Range("M136").Value2 = -2008.783799999999
$ws.Range("N136").Value2 = -10317.75

# BSM (sheet3) row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 1793.5
$ws.Range("I99").Value2 = 1024.875
$ws.Range("J99").Value2 = 2135.111
$ws.Range("K99").Value2 = 1024.875
$ws.Range("L99").Value2 = 2135.111
$ws.Range("M99").Value2 = 473.125
$ws.Range("N99").Value2 = -5131.111

# BSM (sheet3) row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 15152330
$ws.Range("I107").Value2 = 19231408
$ws.Range("K107").Value2 = 19231408
$ws.Range("M107").Value2 = -19229488

# BSM (sheet3) row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value2 = 39960.715
$ws.Range("I135").Value2 = 20000
$ws.Range("J135").Value2 = 47945
$ws.Range("K135").Value2 = 20000
$ws.Range("L135").Value2 = 47945
$ws.Range("N135").Value2 = -58085
$ws.Range("M135").Value2 = -14930

# CRP (sheet4) row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1866.8658
$ws.Range("I31").Value2 = 1066.0927
$ws.Range("J31").Value2 = 3411.2144
$ws.Range("K31").Value2 = 1066.0927
$ws.Range("L31").Value2 = 3411.2144
$ws.Range("M31").Value2 = -771.0926999999999
$ws.Range("N31").Value2 = -4001.2144

# CRP (sheet4) row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 1866.8658
$ws.Range("I34").Value2 = 1066.0927
$ws.Range("J34").Value2 = 3411.2144
$ws.Range("K34").Value2 = 1066.0927
$ws.Range("L34").Value2 = 3411.2144
$ws.Range("M34").Value2 = -864.0926999999999
$ws.Range("N34").Value2 = -3815.2144

# CRP (sheet4) row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 111124850
$ws.Range("J99").Value2 = 3553.5
$ws.Range("L99").Value2 = 3553.5
$ws.Range("N99").Value2 = -6549.5

# CRP (sheet4) row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value2 = 53572304
$ws.Range("I122").Value2 = 150000580
$ws.Range("J122").Value2 = 1043.7778
$ws.Range("K122").Value2 = 450001740
$ws.Range("L122").Value2 = 3131.3334
$ws.Range("M122").Value2 = -449999290
$ws.Range("N122").Value2 = -8031.3334

# CRP (sheet4) row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value2 = 111124850
$ws.Range("J126").Value2 = 3553.5
$ws.Range("L126").Value2 = 10660.5
$ws.Range("N126").Value2 = -15600.5

# CUL (sheet5) row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value2 = 62516960
$ws.Range("I44").Value2 = 100002200
$ws.Range("J44").Value2 = 41563.668
$ws.Range("K44").Value2 = 300006600
$ws.Range("L44").Value2 = 124691.004
$ws.Range("M44").Value2 = -300006202
$ws.Range("N44").Value2 = -125487.004

# CUL (sheet5) row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 55561916
$ws.Range("J68").Value2 = 8602.154
$ws.Range("L68").Value2 = 25806.462
$ws.Range("N68").Value2 = -27428.462

# CUL (sheet5) row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value2 = 17464.428
$ws.Range("J69").Value2 = 24067.8
$ws.Range("L69").Value2 = 72203.39999999999
$ws.Range("N69").Value2 = -73825.39999999999

# CUL (sheet5) row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value2 = 55561916
$ws.Range("J71").Value2 = 8602.154
$ws.Range("L71").Value2 = 77419.386
$ws.Range("N71").Value2 = -85531.386

# CUL (sheet5) row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value2 = 17464.428
$ws.Range("J72").Value2 = 24067.8
$ws.Range("L72").Value2 = 216610.2
$ws.Range("N72").Value2 = -224722.2

# CUL (sheet5) row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 13097.692
$ws.Range("I80").Value2 = 1991
$ws.Range("J80").Value2 = 15117.091
$ws.Range("K80").Value2 = 5973
$ws.Range("L80").Value2 = 45351.273
$ws.Range("M80").Value2 = -5037
$ws.Range("N80").Value2 = -47223.273

# CUL (sheet5) row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value2 = 13097.692
$ws.Range("I83").Value2 = 1991
$ws.Range("J83").Value2 = 15117.091
$ws.Range("K83").Value2 = 17919
$ws.Range("L83").Value2 = 136053.819
$ws.Range("M83").Value2 = -13239
$ws.Range("N83").Value2 = -145413.819

# CUL (sheet5) row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value2 = 1522
$ws.Range("I126").Value2 = 1159.8
$ws.Range("K126").Value2 = 3479.4
$ws.Range("M126").Value2 = 1460.6

# GSM (sheet6) row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 5252.5

# LTW (sheet7) row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value2 = 32929.332
$ws.Range("J108").Value2 = 32929.332
$ws.Range("L108").Value2 = 32929.332
$ws.Range("N108").Value2 = -40609.332

# LTW (sheet7) row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 8409
$ws.Range("I122").Value2 = 9298.964
$ws.Range("J122").Value2 = 4255.8335
$ws.Range("K122").Value2 = 27896.892
$ws.Range("L122").Value2 = 12767.5005
$ws.Range("M122").Value2 = -25446.892
$ws.Range("N122").Value2 = -17667.5005

# LTW (sheet7) row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 4963.426
$ws.Range("I132").Value2 = 1894.4166
$ws.Range("J132").Value2 = 7418.6333
$ws.Range("K132").Value2 = 5683.2498
$ws.Range("L132").Value2 = 22255.8999
$ws.Range("M132").Value2 = -3153.2498
$ws.Range("N132").Value2 = -27315.8999

# WVR (sheet8) row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 2002.2174
$ws.Range("I122").Value2 = 1345.4286
$ws.Range("J122").Value2 = 3023.889
$ws.Range("K122").Value2 = 4036.2858
$ws.Range("L122").Value2 = 9071.667000000001
$ws.Range("M122").Value2 = -1586.2858
$ws.Range("N122").Value2 = -13971.667
